$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.506.23"
$ws.Range("E2").Value = "  +5.61%  "
$ws.Range("D3").Value = "1.722.92"
$ws.Range("E3").Value = "  +4.60%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "225.81"
$ws.Range("E5").Value = "  +3.66%  "
$ws.Range("D6").Value = "0.5345"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  +6.72%  "
$ws.Range("D11").Value = "0.07716"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "1.720.52"
$ws.Range("E13").Value = "  +4.97%  "
$ws.Range("D14").Value = "1.960.97"
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("D15").Value = "0.5831"
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").Value = "0.0₅8286"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "67.99"
$ws.Range("E17").Value = "  +4.59%  "
$ws.Range("D18").Value = "27.506.00"
$ws.Range("E18").Value = "  +5.62%  "
$ws.Range("D19").Value = "219.64"
$ws.Range("E19").Value = "  +15.22%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "4.737"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").Value = "6.079"
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "145.97"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").Value = "1.729"
$ws.Range("E26").Value = "  +14.45%  "
$ws.Range("D27").Value = "0.1236"
$ws.Range("E27").Value = "  +4.84%  "
$ws.Range("D28").Value = "7.402"
$ws.Range("E28").Value = "  +3.42%  "
$ws.Range("D29").Value = "16.56"
$ws.Range("E29").Value = "  +4.96%  "
$ws.Range("D30").Value = "0.05531"
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("E31").Value = "  +3.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.560"
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("D33").Value = "3.446"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("D34").Value = "1.665"
$ws.Range("E34").Value = "  +7.99%  "
$ws.Range("D35").Value = "2.857"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").Value = "0.9655"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").Value = "2.419"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "0.5982"
$ws.Range("E38").Value = "  +7.26%  "
$ws.Range("D39").Value = "0.01653"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.910"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("D42").Value = "1.054.84"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "1.003"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "101.31"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "1.866.96"
$ws.Range("E45").Value = "  +4.47%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").Value = "58.92"
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.230"
$ws.Range("E48").Value = "  +4.35%  "
$ws.Range("D49").Value = "0.4448"
$ws.Range("E49").Value = "  +3.13%  "
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").Value = "0.05237"
$ws.Range("E51").Value = "  +2.05%  "
